$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new category / radio-button pairs in columns C and D for
# rows 3-5 (row 2 already holds "Clinical" / "Clinical_radio_button").
$ws.Range("C3").Value = "Economic"
$ws.Range("D3").Value = "Economic_radio_button"

$ws.Range("C4").Value = "Quality of Life"
$ws.Range("D4").Value = "Quality of Life_radio_button"

$ws.Range("C5").Value = "Real-world Evidence"
$ws.Range("D5").Value = "Real-world Evidence_radio_button"

# Update the sheet view: change the current selection to C3:D5 (anchored at C3).
$ws.Range("C3:D5").Select() | Out-Null
